# validLogin.xlsx update:
#  - Phone number cell (A2) becomes a text value "790000023" (quote-prefixed,
#    so it is stored/displayed as text rather than a number).
#  - Password cell (B2) is updated from "Kalemon123@" to "Kalemon1234@".
#  - A hyperlink is attached to the password cell (B2), which also applies
#    the built-in "Hyperlink" cell style (blue/underlined) to it.
#  - The active selection is moved to B2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: store the phone number as text (leading apostrophe forces a text
# value with a quote prefix instead of a numeric value).
$ws.Range("A2").Value = "'790000023"

# B2: new password text.
$ws.Range("B2").Value = "Kalemon1234@"

# Attach a hyperlink to the password cell.
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.kalemon.com")

# Move the selection to B2, matching the saved workbook state.
$ws.Range("B2").Select()
